$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -------------------------------------------------------------------------
# 1. New retailer alias "Sunil Kr." for PRAKASH GENERAL SRINGAR STORE (row 71)
#    (adds a new shared string and a new value in column B of row 71)
# -------------------------------------------------------------------------
$ws.Range("B49").Copy() | Out-Null
$ws.Range("B71").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("B71").Value = "Sunil Kr."

# -------------------------------------------------------------------------
# 2. New daily collection figures for 30-Dec-2020 (column AJ)
#    Each entry below corresponds to a retailer row; the plain-style rows
#    copy formatting from M12, the highlighted-style rows copy from M13.
# -------------------------------------------------------------------------
$plainStyleRows = @(
    @{Row=12; Value=1000},
    @{Row=14; Value=2000},
    @{Row=20; Value=1500},
    @{Row=22; Value=2000},
    @{Row=23; Value=3000},
    @{Row=25; Value=3000},
    @{Row=26; Value=1000},
    @{Row=27; Value=2000},
    @{Row=28; Value=2000},
    @{Row=41; Value=3000},
    @{Row=69; Value=1500}
)

$highlightStyleRows = @(
    @{Row=13; Value=5000},
    @{Row=18; Value=5000},
    @{Row=35; Value=3000},
    @{Row=49; Value=3000},
    @{Row=71; Value=2000}
)

foreach ($item in $plainStyleRows) {
    $target = "AJ" + $item.Row
    $ws.Range("M12").Copy() | Out-Null
    $ws.Range($target).PasteSpecial(-4122) | Out-Null
    $ws.Range($target).Value = $item.Value
}

foreach ($item in $highlightStyleRows) {
    $target = "AJ" + $item.Row
    $ws.Range("M13").Copy() | Out-Null
    $ws.Range($target).PasteSpecial(-4122) | Out-Null
    $ws.Range($target).Value = $item.Value
}

# -------------------------------------------------------------------------
# 3. Refresh the view state - last worked cell is AJ71, scrolled further
#    down the sheet (used to be near the bottom, now mid-sheet).
# -------------------------------------------------------------------------
$ws.Range("AJ71").Select()
